$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1200.6111
$ws.Range("I40").Value = 798.875
$ws.Range("J40").Value = 1522
$ws.Range("K40").Value = 798.875
$ws.Range("L40").Value = 1522
$ws.Range("M40").Value = -623.875
$ws.Range("N40").Value = -1872

# Row 61
$ws.Range("H61").Value = 1129.2858
$ws.Range("I61").Value = 484.16666
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1452.49998
$ws.Range("L61").Value = 15000
$ws.Range("M61").Value = -1280.49998
$ws.Range("N61").Value = -15344

# Row 96
$ws.Range("H96").Value = 14706670
$ws.Range("I96").Value = 27778240
$ws.Range("J96").Value = 1153.875
$ws.Range("K96").Value = 83334720
$ws.Range("L96").Value = 3461.625
$ws.Range("M96").Value = -83333347
$ws.Range("N96").Value = -6207.625

# Row 103
$ws.Range("H103").Value = 57
$ws.Range("I103").Value = 53.64706
$ws.Range("J103").Value = 66.5
$ws.Range("K103").Value = 160.94118
$ws.Range("L103").Value = 199.5
$ws.Range("M103").Value = 425.05882
$ws.Range("N103").Value = -1371.5

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 137
$ws.Range("H137").Value = 41607.92
$ws.Range("I137").Value = 1249.8125
$ws.Range("J137").Value = 113355.664
$ws.Range("K137").Value = 3749.4375
$ws.Range("L137").Value = 340066.992
$ws.Range("M137").Value = -1199.4375
$ws.Range("N137").Value = -345166.992

# Row 138
$ws.Range("H138").Value = 2527.4583
$ws.Range("I138").Value = 547
$ws.Range("J138").Value = 3342.9412
$ws.Range("K138").Value = 1641
$ws.Range("L138").Value = 10028.8236
$ws.Range("M138").Value = 3499
$ws.Range("N138").Value = -20308.8236

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 30194.29
$ws.Range("I32").Value = 37047.207
$ws.Range("K32").Value = 37047.207
$ws.Range("M32").Value = -36760.207

# Row 88
$ws.Range("H88").Value = 48197.91
$ws.Range("I88").Value = 1534.3334
$ws.Range("K88").Value = 1534.3334
$ws.Range("M88").Value = -1128.3334

# Row 91
$ws.Range("H91").Value = 48197.91
$ws.Range("I91").Value = 1534.3334
$ws.Range("K91").Value = 1534.3334
$ws.Range("M91").Value = -130.3334

# Row 132
$ws.Range("H132").Value = 17322.334
$ws.Range("I132").Value = 1619.3889
$ws.Range("J132").Value = 36165.867
$ws.Range("K132").Value = 4858.1667
$ws.Range("L132").Value = 108497.601
$ws.Range("M132").Value = -2328.1667
$ws.Range("N132").Value = -113557.601

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 31174.53
$ws.Range("I86").Value = 46655.547
$ws.Range("K86").Value = 46655.547
$ws.Range("M86").Value = -45532.547

# Row 89
$ws.Range("H89").Value = 31174.53
$ws.Range("I89").Value = 46655.547
$ws.Range("K89").Value = 233277.735
$ws.Range("M89").Value = -227661.735

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 5798.3887
$ws.Range("J94").Value = 7472.5835
$ws.Range("L94").Value = 7472.5835
$ws.Range("N94").Value = -8374.583500000001

# Row 99
$ws.Range("H99").Value = 5505.263
$ws.Range("I99").Value = 3960
$ws.Range("J99").Value = 7222.222
$ws.Range("K99").Value = 3960
$ws.Range("L99").Value = 7222.222
$ws.Range("M99").Value = -2462
$ws.Range("N99").Value = -10218.222

# Row 126
$ws.Range("H126").Value = 5505.263
$ws.Range("I126").Value = 3960
$ws.Range("J126").Value = 7222.222
$ws.Range("K126").Value = 11880
$ws.Range("L126").Value = 21666.666
$ws.Range("M126").Value = -9410
$ws.Range("N126").Value = -26606.666

# Row 134
$ws.Range("H134").Value = 800.92
$ws.Range("I134").Value = 667.4
$ws.Range("J134").Value = 1001.2
$ws.Range("K134").Value = 2002.2
$ws.Range("L134").Value = 3003.6
$ws.Range("M134").Value = 532.8000000000002
$ws.Range("N134").Value = -8073.6

# Row 135
$ws.Range("H135").Value = 50933.332
$ws.Range("J135").Value = 50933.332
$ws.Range("L135").Value = 50933.332
$ws.Range("N135").Value = -61073.332

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1280.3914
$ws.Range("J68").Value = 1415.6842
$ws.Range("L68").Value = 4247.0526
$ws.Range("N68").Value = -5869.0526

# Row 71
$ws.Range("H71").Value = 1280.3914
$ws.Range("J71").Value = 1415.6842
$ws.Range("L71").Value = 12741.1578
$ws.Range("N71").Value = -20853.1578

# Row 107
$ws.Range("H107").Value = 4500.5
$ws.Range("I107").Value = 8823.833000000001
$ws.Range("J107").Value = 794.7857
$ws.Range("K107").Value = 26471.499
$ws.Range("L107").Value = 2384.3571
$ws.Range("M107").Value = -24551.499
$ws.Range("N107").Value = -6224.3571

# Row 109
$ws.Range("H109").Value = 5836.75
$ws.Range("J109").Value = 6284.857
$ws.Range("L109").Value = 18854.571
$ws.Range("N109").Value = -20934.571

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3231.5833
$ws.Range("J113").Value = 4800
$ws.Range("L113").Value = 4800
$ws.Range("N113").Value = -9140

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1068.55
$ws.Range("I16").Value = 842.4375
$ws.Range("J16").Value = 1973
$ws.Range("K16").Value = 842.4375
$ws.Range("L16").Value = 1973
$ws.Range("M16").Value = -672.4375
$ws.Range("N16").Value = -2313

# Row 22
$ws.Range("H22").Value = 2316.2
$ws.Range("I22").Value = 2795.25
$ws.Range("K22").Value = 2795.25
$ws.Range("M22").Value = -2500.25

# Row 27
$ws.Range("H27").Value = 2316.2
$ws.Range("I27").Value = 2795.25
$ws.Range("K27").Value = 2795.25
$ws.Range("M27").Value = -2688.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1376.625
$ws.Range("I81").Value = 1316.1428
$ws.Range("K81").Value = 2632.2856
$ws.Range("M81").Value = -1571.2856

# Row 84
$ws.Range("H84").Value = 1376.625
$ws.Range("I84").Value = 1316.1428
$ws.Range("K84").Value = 13161.428
$ws.Range("M84").Value = -7857.428

# Row 111
$ws.Range("H111").Value = 33000
$ws.Range("J111").Value = 33000
$ws.Range("L111").Value = 33000
$ws.Range("N111").Value = -41180

# Row 132
$ws.Range("H132").Value = 1367.2391
$ws.Range("I132").Value = 1054.5588
$ws.Range("J132").Value = 2253.1667
$ws.Range("K132").Value = 3163.6764
$ws.Range("L132").Value = 6759.500100000001
$ws.Range("M132").Value = -633.6764000000003
$ws.Range("N132").Value = -11819.5001
